$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values must stay text (they use a dot as a
# thousands separator, e.g. "68.176.04", or carry significant
# trailing zeros, e.g. "1.00"), so force a text number format on
# each such cell before writing its value, then restore the default
# (unstyled) style so no stray formatting is left behind.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "68.176.04"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.59%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.800.08"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.14%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "601.16"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "165.14"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.24%  "
$ws.Range("E7").Value = "  -0.01%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.518"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("E10").Value = "  +0.57%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.46"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("E12").Value = "  -1.36%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "35.85"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.66%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.438.02"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.800.71"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.16%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "68.198.21"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "18.45"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("E19").Value = "  -0.16%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "463.19"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.40%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "9.72"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.13%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.701"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.0000149"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.31%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "83.01"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  +0.09%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.01"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.69%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "3.949.29"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  -4.91%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.34"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.31%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.23"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  -0.49%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.0996"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.139"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.31"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("E39").Value = "  +0.67%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.987"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  +0.00%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.301"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.53%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "47.46"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.39%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "43.16"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.50%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "151.79"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  +2.36%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "395.90"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.36"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +5.71%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "26.77"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.02%  "
